$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row labels: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410"
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v -like "*_old") {
            $cell.Value2 = ($v -replace "_old$", "_FV2404")
        } elseif ($v -like "*_new") {
            $cell.Value2 = ($v -replace "_new$", "_FV2410")
        }
    }
}

# 2. Turn the used range into an Excel Table ("Table1") with the renamed headers
$rng = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# 3. Freeze the header row (split pane under row 1)
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
